$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the new, longer values added below
$ws.Columns.Item(1).ColumnWidth = 16

# New "other in CSV" section header (bold, like the existing section headers)
$ws.Range("A29").Value = "other in CSV"
$ws.Range("A29").Font.Bold = $true

# New row describing the SOIL_STONES field
$ws.Range("A30").Value = "SOIL_STONES"
$ws.Range("B30").Value = "Percentage of Stones in the sample, as a class"
$ws.Range("D30").Value = "?"

# Move the active selection to B13, matching the saved workbook state
$ws.Range("B13").Select()
